$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 642
$ws.Range("I2").Value = 1701
$ws.Range("J2").Value = 7111
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 1956
$ws.Range("M2").Value = 134
$ws.Range("N2").Value = 1254
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 38
$ws.Range("Q2").Value = 12
$ws.Range("R2").Value = 83
$ws.Range("S2").Value = 779
$ws.Range("T2").Value = 1214
$ws.Range("U2").Value = 97
$ws.Range("V2").Value = 11171
$ws.Range("X2").Value = 11499
$ws.Range("Y2").Value = 20
$ws.Range("Z2").Value = 135
$ws.Range("AA2").Value = 74
